$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Characters and strings) - "What you learn" and "Estimated Time"
$ws.Range("D9").Value = "String build in functions"
$ws.Range("E9").Value = "2 weeks"

# Row 8 (Pointers) - "Question number" and row 9 "Question number"
$ws.Range("G8").Value = "6 problems"
$ws.Range("G9").Value = "WIP"

# Row 8 Comments - "Bad" style
$ws.Range("H8").Value = "Card game, Maze traveller are remaining."
$ws.Range("H8").Style = "Bad"
$ws.Range("H8").HorizontalAlignment = -4131
$ws.Range("H8").VerticalAlignment = -4108

# Row 10 (Formatted input/output) - "Neutral" style
$ws.Range("F10").Value = " No exercise"
$ws.Range("F10").Style = "Neutral"
$ws.Range("F10").HorizontalAlignment = -4108
$ws.Range("F10").VerticalAlignment = -4130

# Row 11 (Structure, union, Bit manipulation) - "Neutral" style
$ws.Range("F11").Value = "No exercise"
$ws.Range("F11").Style = "Neutral"
$ws.Range("F11").HorizontalAlignment = -4108
$ws.Range("F11").VerticalAlignment = -4108

# Row 12 (File processing) - "Neutral" style
$ws.Range("F12").Value = "No exercise"
$ws.Range("F12").Style = "Neutral"
$ws.Range("F12").HorizontalAlignment = -4108
$ws.Range("F12").VerticalAlignment = -4108

# Selection / view change
$ws.Range("F14").Select()
